# Only measure encoding and solving time of SAT/SCL approaches & rerun.
# Splits the old single "wall_time_sec" column into separate
# encoding_time_sec / solving_time_sec / total_time_sec columns,
# shifts peak_memory_mb / total_clauses / total_variables right,
# and fills in the new timing figures from the rerun experiment.

$wb = $excel.ActiveWorkbook

#### Sheet: n1_d40 ####
$ws = $wb.Worksheets.Item("n1_d40")

$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"
$ws.Range("I1").Value = "peak_memory_mb"
$ws.Range("J1").Value = "total_clauses"
$ws.Range("K1").Value = "total_variables"

$ws.Range("F2").Value = 0.0112384
$ws.Range("G2").Value = 0.00429377
$ws.Range("H2").Value = 0.0155321
$ws.Range("I2").Value = 5.7
$ws.Range("F3").Value = 0.0112592
$ws.Range("G3").Value = 0.00439118
$ws.Range("H3").Value = 0.0156504
$ws.Range("I3").Value = 5.7
$ws.Range("F4").Value = 0.0114624
$ws.Range("G4").Value = 0.00429688
$ws.Range("H4").Value = 0.0157593
$ws.Range("I4").Value = 5.7
$ws.Range("F5").Value = 0.0117934
$ws.Range("G5").Value = 0.00436984
$ws.Range("H5").Value = 0.0161632
$ws.Range("I5").Value = 5.6
$ws.Range("F6").Value = 0.0113563
$ws.Range("G6").Value = 0.00441622
$ws.Range("H6").Value = 0.0157726
$ws.Range("I6").Value = 5.7
$ws.Range("F7").Value = 0.0111856
$ws.Range("G7").Value = 0.00431165
$ws.Range("H7").Value = 0.0154972
$ws.Range("I7").Value = 5.6
$ws.Range("F8").Value = 0.0113704
$ws.Range("G8").Value = 0.00426493
$ws.Range("H8").Value = 0.0156353
$ws.Range("I8").Value = 5.6
$ws.Range("F9").Value = 0.0112725
$ws.Range("G9").Value = 0.00425964
$ws.Range("H9").Value = 0.0155321
$ws.Range("I9").Value = 5.7
$ws.Range("F10").Value = 0.0112101
$ws.Range("G10").Value = 0.0042701
$ws.Range("H10").Value = 0.0154802
$ws.Range("I10").Value = 5.6
$ws.Range("F11").Value = 0.0110138
$ws.Range("G11").Value = 0.00425581
$ws.Range("H11").Value = 0.0152696
$ws.Range("I11").Value = 5.6
$ws.Range("F12").Value = 0.01131621
$ws.Range("G12").Value = 0.004313002
$ws.Range("H12").Value = 0.0156292
$ws.Range("I12").Value = 5.65

#### Sheet: n1_d60 ####
$ws = $wb.Worksheets.Item("n1_d60")

$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"
$ws.Range("I1").Value = "peak_memory_mb"
$ws.Range("J1").Value = "total_clauses"
$ws.Range("K1").Value = "total_variables"

$ws.Range("F2").Value = 0.0202972
$ws.Range("G2").Value = 0.00517888
$ws.Range("H2").Value = 0.0254761
$ws.Range("I2").Value = 8.300000000000001
$ws.Range("F3").Value = 0.0199553
$ws.Range("G3").Value = 0.00513706
$ws.Range("H3").Value = 0.0250924
$ws.Range("I3").Value = 8.6
$ws.Range("F4").Value = 0.0203587
$ws.Range("G4").Value = 0.00506974
$ws.Range("H4").Value = 0.0254285
$ws.Range("I4").Value = 8.199999999999999
$ws.Range("F5").Value = 0.0196807
$ws.Range("G5").Value = 0.00486002
$ws.Range("H5").Value = 0.0245407
$ws.Range("I5").Value = 8.699999999999999
$ws.Range("F6").Value = 0.0198155
$ws.Range("G6").Value = 0.00478992
$ws.Range("H6").Value = 0.0246054
$ws.Range("I6").Value = 8.6
$ws.Range("F7").Value = 0.0194543
$ws.Range("G7").Value = 0.00478416
$ws.Range("H7").Value = 0.0242385
$ws.Range("I7").Value = 8.699999999999999
$ws.Range("F8").Value = 0.0195151
$ws.Range("G8").Value = 0.00491374
$ws.Range("H8").Value = 0.0244289
$ws.Range("I8").Value = 8.699999999999999
$ws.Range("F9").Value = 0.0200908
$ws.Range("G9").Value = 0.004941
$ws.Range("H9").Value = 0.0250318
$ws.Range("I9").Value = 8.699999999999999
$ws.Range("F10").Value = 0.0200334
$ws.Range("G10").Value = 0.00491252
$ws.Range("H10").Value = 0.0249459
$ws.Range("I10").Value = 8.6
$ws.Range("F11").Value = 0.0200459
$ws.Range("G11").Value = 0.00501891
$ws.Range("H11").Value = 0.0250648
$ws.Range("I11").Value = 8.6
$ws.Range("F12").Value = 0.01992469
$ws.Range("G12").Value = 0.004960595
$ws.Range("H12").Value = 0.0248853
$ws.Range("I12").Value = 8.569999999999999

#### Sheet: n1_d80 ####
$ws = $wb.Worksheets.Item("n1_d80")

$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"
$ws.Range("I1").Value = "peak_memory_mb"
$ws.Range("J1").Value = "total_clauses"
$ws.Range("K1").Value = "total_variables"

$ws.Range("F2").Value = 0.0276765
$ws.Range("G2").Value = 0.00511364
$ws.Range("H2").Value = 0.0327902
$ws.Range("I2").Value = 11.5
$ws.Range("F3").Value = 0.0273391
$ws.Range("G3").Value = 0.00514643
$ws.Range("H3").Value = 0.0324856
$ws.Range("I3").Value = 11.5
$ws.Range("F4").Value = 0.0275934
$ws.Range("G4").Value = 0.00510339
$ws.Range("H4").Value = 0.0326968
$ws.Range("I4").Value = 11.4
$ws.Range("F5").Value = 0.0281174
$ws.Range("G5").Value = 0.00519673
$ws.Range("H5").Value = 0.0333142
$ws.Range("I5").Value = 11.5
$ws.Range("F6").Value = 0.0275159
$ws.Range("G6").Value = 0.00503877
$ws.Range("H6").Value = 0.0325547
$ws.Range("I6").Value = 11.5
$ws.Range("F7").Value = 0.0282104
$ws.Range("G7").Value = 0.00524271
$ws.Range("H7").Value = 0.0334531
$ws.Range("I7").Value = 11.4
$ws.Range("F8").Value = 0.0273908
$ws.Range("G8").Value = 0.00508815
$ws.Range("H8").Value = 0.032479
$ws.Range("I8").Value = 11.4
$ws.Range("F9").Value = 0.0278732
$ws.Range("G9").Value = 0.00526387
$ws.Range("H9").Value = 0.0331371
$ws.Range("I9").Value = 11.6
$ws.Range("F10").Value = 0.0275159
$ws.Range("G10").Value = 0.00506259
$ws.Range("H10").Value = 0.0325785
$ws.Range("I10").Value = 11.5
$ws.Range("F11").Value = 0.0275894
$ws.Range("G11").Value = 0.00506641
$ws.Range("H11").Value = 0.0326558
$ws.Range("I11").Value = 11.6
$ws.Range("F12").Value = 0.0276822
$ws.Range("G12").Value = 0.005132269
$ws.Range("H12").Value = 0.0328145
$ws.Range("I12").Value = 11.49

#### Sheet: n1_d100 ####
$ws = $wb.Worksheets.Item("n1_d100")

$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"
$ws.Range("I1").Value = "peak_memory_mb"
$ws.Range("J1").Value = "total_clauses"
$ws.Range("K1").Value = "total_variables"

$ws.Range("F2").Value = 0.036243
$ws.Range("G2").Value = 0.00527944
$ws.Range("H2").Value = 0.0415225
$ws.Range("I2").Value = 13.4
$ws.Range("F3").Value = 0.0368803
$ws.Range("G3").Value = 0.00543953
$ws.Range("H3").Value = 0.0423199
$ws.Range("I3").Value = 13.3
$ws.Range("F4").Value = 0.0362934
$ws.Range("G4").Value = 0.00580882
$ws.Range("H4").Value = 0.0421022
$ws.Range("I4").Value = 13.3
$ws.Range("F5").Value = 0.0367109
$ws.Range("G5").Value = 0.00536601
$ws.Range("H5").Value = 0.0420769
$ws.Range("I5").Value = 13.4
$ws.Range("F6").Value = 0.0376996
$ws.Range("G6").Value = 0.00590281
$ws.Range("H6").Value = 0.0436024
$ws.Range("I6").Value = 13.4
$ws.Range("F7").Value = 0.0367168
$ws.Range("G7").Value = 0.00584487
$ws.Range("H7").Value = 0.0425617
$ws.Range("I7").Value = 13.3
$ws.Range("F8").Value = 0.0372151
$ws.Range("G8").Value = 0.00635291
$ws.Range("H8").Value = 0.043568
$ws.Range("I8").Value = 13.3
$ws.Range("F9").Value = 0.036205
$ws.Range("G9").Value = 0.00554181
$ws.Range("H9").Value = 0.0417468
$ws.Range("I9").Value = 13.3
$ws.Range("F10").Value = 0.0367507
$ws.Range("G10").Value = 0.00565785
$ws.Range("H10").Value = 0.0424085
$ws.Range("I10").Value = 13.4
$ws.Range("F11").Value = 0.0362606
$ws.Range("G11").Value = 0.00536924
$ws.Range("H11").Value = 0.0416299
$ws.Range("I11").Value = 13.4
$ws.Range("F12").Value = 0.03669753999999999
$ws.Range("G12").Value = 0.005656329
$ws.Range("H12").Value = 0.04235388
$ws.Range("I12").Value = 13.35
